$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values per the diff (shared strings changed, B2 numeric value changed)
$ws.Range("A2").Value = "lSqbS113"
$ws.Range("B2").Value = 23103015
$ws.Range("C2").Value = "lvocunc96"
$ws.Range("D2").Value = "h4k%9WE&"
$ws.Range("F2").Value = "JNThQGbk"
$ws.Range("G2").Value = "fCxn"
